# Regenerate the "K" column (column G) values in the save_data sheet.
# These values were recomputed (K instead of the old Strike# derived value,
# std/mean and s_vals recalculated upstream) and are written back here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Value = 0 },
    @{ Row = 3;  Value = 0 },
    @{ Row = 4;  Value = 0 },
    @{ Row = 5;  Value = 1 },
    @{ Row = 6;  Value = 0 },
    @{ Row = 7;  Value = 1 },
    @{ Row = 8;  Value = 1 },
    @{ Row = 10; Value = 1 },
    @{ Row = 11; Value = 0 },
    @{ Row = 12; Value = 1 },
    @{ Row = 13; Value = 0 },
    @{ Row = 14; Value = 0 },
    @{ Row = 15; Value = 0 },
    @{ Row = 16; Value = 0 },
    @{ Row = 17; Value = 2 },
    @{ Row = 18; Value = 1 },
    @{ Row = 19; Value = 1 },
    @{ Row = 20; Value = 1 },
    @{ Row = 21; Value = 0 },
    @{ Row = 22; Value = 1 },
    @{ Row = 23; Value = 3 },
    @{ Row = 24; Value = 1 },
    @{ Row = 25; Value = 0 },
    @{ Row = 26; Value = 0 },
    @{ Row = 27; Value = 0 },
    @{ Row = 28; Value = 1 },
    @{ Row = 29; Value = 0 },
    @{ Row = 30; Value = 1 },
    @{ Row = 31; Value = 0 },
    @{ Row = 32; Value = 1 },
    @{ Row = 33; Value = 0 },
    @{ Row = 34; Value = 2 },
    @{ Row = 35; Value = 0 },
    @{ Row = 36; Value = 1 },
    @{ Row = 37; Value = 0 },
    @{ Row = 38; Value = 1 },
    @{ Row = 39; Value = 0 },
    @{ Row = 40; Value = 1 },
    @{ Row = 41; Value = 0 },
    @{ Row = 42; Value = 0 },
    @{ Row = 43; Value = 2 },
    @{ Row = 45; Value = 1 },
    @{ Row = 46; Value = 0 },
    @{ Row = 48; Value = 1 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}
